$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)
$rows = @(
    @{ Row = 2; D = 45028; J = 1000; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 3; D = 45020; J = 1200; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 4; D = 45077; J = 760; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 5; D = 44978; J = 1000; K = 1800; L = 2000; M = 1900; P = 633 }
    @{ Row = 6; D = 44970; J = 800; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 7; D = 45084; J = 900; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 8; D = 45041; J = 1160; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 9; D = 45013; J = 1100; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 10; D = 44964; J = 1000; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 11; D = 45034; J = 1100; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 12; D = 44951; J = 800; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 13; D = 45070; J = 800; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 14; D = 45091; J = 800; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 15; D = 44911; J = 700; K = 1800; L = 2000; M = 1900; P = 633 }
    @{ Row = 16; D = 44827; J = 1200; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 17; D = 44985; J = 1000; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 18; D = 44999; J = 1100; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 19; D = 44965; J = 1120; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 20; D = 44848; J = 1000; K = 1500; L = 2000; M = 1750; P = 583 }
    @{ Row = 21; D = 45062; J = 1100; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 22; D = 44685; J = 400; K = 1500; L = 2000; M = 1750; P = 583 }
    @{ Row = 23; D = 45035; J = 1100; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 24; D = 44992; J = 1040; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 25; D = 44910; J = 1000; K = 1800; L = 2000; M = 1900; P = 633 }
    @{ Row = 26; D = 45006; J = 1100; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 27; D = 44953; J = 1000; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 28; D = 44883; J = 500; K = 1800; L = 2000; M = 1900; P = 633 }
    @{ Row = 29; D = 45007; J = 1160; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 30; D = 44971; J = 1000; K = 2000; L = 2500; M = 2250; P = 750 }
    @{ Row = 31; D = 44881; J = 500; K = 1900; L = 2000; M = 1950; P = 650 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value  = $r.D   # Fecha
    $ws.Cells.Item($r.Row, 10).Value = $r.J   # Volumen
    $ws.Cells.Item($r.Row, 11).Value = $r.K   # Precio minimo
    $ws.Cells.Item($r.Row, 12).Value = $r.L   # Precio maximo
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # Precio promedio ponderado
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # Precio $/Kg
}
